$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 82.73729466666667
$ws.Range("H2").Value = 248.211884
$ws.Range("I2").Value = 0.09847102321391109
$ws.Range("J2").Value = 0.09847102321391106
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.744736333333333
$ws.Range("N2").Value = 20.234209
$ws.Range("O2").Value = 0.01049273071342527
$ws.Range("P2").Value = 0.01049273071342527
$ws.Range("Q2").Value = 558.0412374599729
$ws.Range("R2").Value = 5022.371137139756
$ws.Range("S2").Value = 0.001033229929659018
$ws.Range("T2").Value = 0.001033229929659018

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 82.73729466666667
$ws.Range("H3").Value = 248.211884
$ws.Range("I3").Value = 0.09847102321391109
$ws.Range("J3").Value = 0.09847102321391106
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 182.5316823333334
$ws.Range("N3").Value = 547.595047
$ws.Range("O3").Value = 0.283963033503136
$ws.Range("P3").Value = 0.2839630335031361
$ws.Range("Q3").Value = 15102.1775872154
$ws.Range("R3").Value = 135919.5982849385
$ws.Range("S3").Value = 0.02796213046397992
$ws.Range("T3").Value = 0.02796213046397992

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 82.73729466666667
$ws.Range("H4").Value = 248.211884
$ws.Range("I4").Value = 0.09847102321391109
$ws.Range("J4").Value = 0.09847102321391106
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 249.1329143333333
$ws.Range("N4").Value = 747.398743
$ws.Range("O4").Value = 0.3875740210972192
$ws.Range("P4").Value = 0.3875740210972192
$ws.Range("Q4").Value = 20612.58334436243
$ws.Range("R4").Value = 185513.2500992618
$ws.Range("S4").Value = 0.03816481042857314
$ws.Range("T4").Value = 0.03816481042857313

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 82.73729466666667
$ws.Range("H5").Value = 248.211884
$ws.Range("I5").Value = 0.09847102321391109
$ws.Range("J5").Value = 0.09847102321391106
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 186.9310863333334
$ws.Range("N5").Value = 560.793259
$ws.Range("O5").Value = 0.2908071500393791
$ws.Range("P5").Value = 0.2908071500393791
$ws.Range("Q5").Value = 15466.17237232111
$ws.Range("R5").Value = 139195.55135089
$ws.Range("S5").Value = 0.02863607762229902
$ws.Range("T5").Value = 0.02863607762229901

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 82.73729466666667
$ws.Range("H6").Value = 248.211884
$ws.Range("I6").Value = 0.09847102321391109
$ws.Range("J6").Value = 0.09847102321391106
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.46044133333334
$ws.Range("N6").Value = 52.38132400000001
$ws.Range("O6").Value = 0.02716306464684043
$ws.Range("P6").Value = 0.02716306464684043
$ws.Range("Q6").Value = 1444.629679606046
$ws.Range("R6").Value = 13001.66711645442
$ws.Range("S6").Value = 0.002674774769399992
$ws.Range("T6").Value = 0.002674774769399991

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 237.0718893333334
$ws.Range("H7").Value = 711.215668
$ws.Range("I7").Value = 0.2821546391135941
$ws.Range("J7").Value = 0.2821546391135941
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.744736333333333
$ws.Range("N7").Value = 20.234209
$ws.Range("O7").Value = 0.01049273071342527
$ws.Range("P7").Value = 0.01049273071342527
$ws.Range("Q7").Value = 1598.987385598512
$ws.Range("R7").Value = 14390.88647038661
$ws.Range("S7").Value = 0.002960572647762633
$ws.Range("T7").Value = 0.002960572647762632

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 237.0718893333334
$ws.Range("H8").Value = 711.215668
$ws.Range("I8").Value = 0.2821546391135941
$ws.Range("J8").Value = 0.2821546391135941
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 182.5316823333334
$ws.Range("N8").Value = 547.595047
$ws.Range("O8").Value = 0.283963033503136
$ws.Range("P8").Value = 0.2839630335031361
$ws.Range("Q8").Value = 43273.13079395516
$ws.Range("R8").Value = 389458.1771455964
$ws.Range("S8").Value = 0.08012148723967878
$ws.Range("T8").Value = 0.08012148723967878

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 237.0718893333334
$ws.Range("H9").Value = 711.215668
$ws.Range("I9").Value = 0.2821546391135941
$ws.Range("J9").Value = 0.2821546391135941
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 249.1329143333333
$ws.Range("N9").Value = 747.398743
$ws.Range("O9").Value = 0.3875740210972192
$ws.Range("P9").Value = 0.3875740210972192
$ws.Range("Q9").Value = 59062.41069612282
$ws.Range("R9").Value = 531561.6962651054
$ws.Range("S9").Value = 0.1093558080524904
$ws.Range("T9").Value = 0.1093558080524904

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 237.0718893333334
$ws.Range("H10").Value = 711.215668
$ws.Range("I10").Value = 0.2821546391135941
$ws.Range("J10").Value = 0.2821546391135941
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 186.9310863333334
$ws.Range("N10").Value = 560.793259
$ws.Range("O10").Value = 0.2908071500393791
$ws.Range("P10").Value = 0.2908071500393791
$ws.Range("Q10").Value = 44316.10581217579
$ws.Range("R10").Value = 398844.9523095821
$ws.Range("S10").Value = 0.08205258647101382
$ws.Range("T10").Value = 0.08205258647101381

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 237.0718893333334
$ws.Range("H11").Value = 711.215668
$ws.Range("I11").Value = 0.2821546391135941
$ws.Range("J11").Value = 0.2821546391135941
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 17.46044133333334
$ws.Range("N11").Value = 52.38132400000001
$ws.Range("O11").Value = 0.02716306464684043
$ws.Range("P11").Value = 0.02716306464684043
$ws.Range("Q11").Value = 4139.37981548716
$ws.Range("R11").Value = 37254.41833938444
$ws.Range("S11").Value = 0.007664184702648488
$ws.Range("T11").Value = 0.007664184702648487

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 277.340215
$ws.Range("H12").Value = 832.0206450000001
$ws.Range("I12").Value = 0.3300805865050695
$ws.Range("J12").Value = 0.3300805865050694
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.744736333333333
$ws.Range("N12").Value = 20.234209
$ws.Range("O12").Value = 0.01049273071342527
$ws.Range("P12").Value = 0.01049273071342527
$ws.Range("Q12").Value = 1870.586624804978
$ws.Range("R12").Value = 16835.27962324481
$ws.Range("S12").Value = 0.00346344670792717
$ws.Range("T12").Value = 0.00346344670792717

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 277.340215
$ws.Range("H13").Value = 832.0206450000001
$ws.Range("I13").Value = 0.3300805865050695
$ws.Range("J13").Value = 0.3300805865050694
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 182.5316823333334
$ws.Range("N13").Value = 547.595047
$ws.Range("O13").Value = 0.283963033503136
$ws.Range("P13").Value = 0.2839630335031361
$ws.Range("Q13").Value = 50623.37602263837
$ws.Range("R13").Value = 455610.3842037453
$ws.Range("S13").Value = 0.09373068464447383
$ws.Range("T13").Value = 0.09373068464447383

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 277.340215
$ws.Range("H14").Value = 832.0206450000001
$ws.Range("I14").Value = 0.3300805865050695
$ws.Range("J14").Value = 0.3300805865050694
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 249.1329143333333
$ws.Range("N14").Value = 747.398743
$ws.Range("O14").Value = 0.3875740210972192
$ws.Range("P14").Value = 0.3875740210972192
$ws.Range("Q14").Value = 69094.57602478325
$ws.Range("R14").Value = 621851.1842230492
$ws.Range("S14").Value = 0.1279306601978983
$ws.Range("T14").Value = 0.1279306601978983

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 277.340215
$ws.Range("H15").Value = 832.0206450000001
$ws.Range("I15").Value = 0.3300805865050695
$ws.Range("J15").Value = 0.3300805865050694
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 186.9310863333334
$ws.Range("N15").Value = 560.793259
$ws.Range("O15").Value = 0.2908071500393791
$ws.Range("P15").Value = 0.2908071500393791
$ws.Range("Q15").Value = 51843.50767387023
$ws.Range("R15").Value = 466591.5690648321
$ws.Range("S15").Value = 0.09598979464486598
$ws.Range("T15").Value = 0.09598979464486597

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 277.340215
$ws.Range("H16").Value = 832.0206450000001
$ws.Range("I16").Value = 0.3300805865050695
$ws.Range("J16").Value = 0.3300805865050694
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 17.46044133333334
$ws.Range("N16").Value = 52.38132400000001
$ws.Range("O16").Value = 0.02716306464684043
$ws.Range("P16").Value = 0.02716306464684043
$ws.Range("Q16").Value = 4842.482553381554
$ws.Range("R16").Value = 43582.34298043399
$ws.Range("S16").Value = 0.008966000309904206
$ws.Range("T16").Value = 0.008966000309904206

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 191.058024
$ws.Range("H17").Value = 573.174072
$ws.Range("I17").Value = 0.22739055213619
$ws.Range("J17").Value = 0.22739055213619
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.744736333333333
$ws.Range("N17").Value = 20.234209
$ws.Range("O17").Value = 0.01049273071342527
$ws.Range("P17").Value = 0.01049273071342527
$ws.Range("Q17").Value = 1288.635996247672
$ws.Range("R17").Value = 11597.72396622905
$ws.Range("S17").Value = 0.002385947830342132
$ws.Range("T17").Value = 0.002385947830342132

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 191.058024
$ws.Range("H18").Value = 573.174072
$ws.Range("I18").Value = 0.22739055213619
$ws.Range("J18").Value = 0.22739055213619
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 182.5316823333334
$ws.Range("N18").Value = 547.595047
$ws.Range("O18").Value = 0.283963033503136
$ws.Range("P18").Value = 0.2839630335031361
$ws.Range("Q18").Value = 34874.14254400238
$ws.Range("R18").Value = 313867.2828960214
$ws.Range("S18").Value = 0.06457051097454551
$ws.Range("T18").Value = 0.06457051097454553

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 191.058024
$ws.Range("H19").Value = 573.174072
$ws.Range("I19").Value = 0.22739055213619
$ws.Range("J19").Value = 0.22739055213619
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 249.1329143333333
$ws.Range("N19").Value = 747.398743
$ws.Range("O19").Value = 0.3875740210972192
$ws.Range("P19").Value = 0.3875740210972192
$ws.Range("Q19").Value = 47598.84232588795
$ws.Range("R19").Value = 428389.5809329915
$ws.Range("S19").Value = 0.08813067065094002
$ws.Range("T19").Value = 0.08813067065094002

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 191.058024
$ws.Range("H20").Value = 573.174072
$ws.Range("I20").Value = 0.22739055213619
$ws.Range("J20").Value = 0.22739055213619
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 186.9310863333334
$ws.Range("N20").Value = 560.793259
$ws.Range("O20").Value = 0.2908071500393791
$ws.Range("P20").Value = 0.2908071500393791
$ws.Range("Q20").Value = 35714.68397902008
$ws.Range("R20").Value = 321432.1558111807
$ws.Range("S20").Value = 0.06612679841260626
$ws.Range("T20").Value = 0.06612679841260624

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 191.058024
$ws.Range("H21").Value = 573.174072
$ws.Range("I21").Value = 0.22739055213619
$ws.Range("J21").Value = 0.22739055213619
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 17.46044133333334
$ws.Range("N21").Value = 52.38132400000001
$ws.Range("O21").Value = 0.02716306464684043
$ws.Range("P21").Value = 0.02716306464684043
$ws.Range("Q21").Value = 3335.957419314593
$ws.Range("R21").Value = 30023.61677383133
$ws.Range("S21").Value = 0.006176624267756068
$ws.Range("T21").Value = 0.006176624267756068

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 52.01228800000001
$ws.Range("H22").Value = 156.036864
$ws.Range("I22").Value = 0.0619031990312353
$ws.Range("J22").Value = 0.06190319903123528
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 6.744736333333333
$ws.Range("N22").Value = 20.234209
$ws.Range("O22").Value = 0.01049273071342527
$ws.Range("P22").Value = 0.01049273071342527
$ws.Range("Q22").Value = 350.8091686533974
$ws.Range("R22").Value = 3157.282517880576
$ws.Range("S22").Value = 0.0006495335977343201
$ws.Range("T22").Value = 0.0006495335977343201

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 52.01228800000001
$ws.Range("H23").Value = 156.036864
$ws.Range("I23").Value = 0.0619031990312353
$ws.Range("J23").Value = 0.06190319903123528
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 182.5316823333334
$ws.Range("N23").Value = 547.595047
$ws.Range("O23").Value = 0.283963033503136
$ws.Range("P23").Value = 0.2839630335031361
$ws.Range("Q23").Value = 9493.890430645848
$ws.Range("R23").Value = 85445.01387581261
$ws.Range("S23").Value = 0.01757822018045797
$ws.Range("T23").Value = 0.01757822018045797

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 52.01228800000001
$ws.Range("H24").Value = 156.036864
$ws.Range("I24").Value = 0.0619031990312353
$ws.Range("J24").Value = 0.06190319903123528
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 249.1329143333333
$ws.Range("N24").Value = 747.398743
$ws.Range("O24").Value = 0.3875740210972192
$ws.Range("P24").Value = 0.3875740210972192
$ws.Range("Q24").Value = 12957.97289058466
$ws.Range("R24").Value = 116621.756015262
$ws.Range("S24").Value = 0.02399207176731735
$ws.Range("T24").Value = 0.02399207176731735

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 52.01228800000001
$ws.Range("H25").Value = 156.036864
$ws.Range("I25").Value = 0.0619031990312353
$ws.Range("J25").Value = 0.06190319903123528
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 186.9310863333334
$ws.Range("N25").Value = 560.793259
$ws.Range("O25").Value = 0.2908071500393791
$ws.Range("P25").Value = 0.2908071500393791
$ws.Range("Q25").Value = 9722.7134985222
$ws.Range("R25").Value = 87504.42148669978
$ws.Range("S25").Value = 0.01800189288859399
$ws.Range("T25").Value = 0.01800189288859399

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 52.01228800000001
$ws.Range("H26").Value = 156.036864
$ws.Range("I26").Value = 0.0619031990312353
$ws.Range("J26").Value = 0.06190319903123528
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 17.46044133333334
$ws.Range("N26").Value = 52.38132400000001
$ws.Range("O26").Value = 0.02716306464684043
$ws.Range("P26").Value = 0.02716306464684043
$ws.Range("Q26").Value = 908.1575032364375
$ws.Range("R26").Value = 8173.417529127937
$ws.Range("S26").Value = 0.001681480597131674
$ws.Range("T26").Value = 0.001681480597131674
